$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 8,31

$arr[0,0] = -2.978358019635732
$arr[0,1] = 1.46110522798494
$arr[0,2] = 1.76488022812725
$arr[0,3] = -0.1819887215705478
$arr[0,4] = -0.0649079407850771
$arr[0,5] = -0.0007307741227984005
$arr[0,6] = -7.557823177602677
$arr[0,7] = 4.290037286492535
$arr[0,8] = 5.37823968090039
$arr[0,9] = -1.613803115652376
$arr[0,10] = -0.4906743893445039
$arr[0,11] = -0.005976284788717497
$arr[0,12] = -4.882783070365717
$arr[0,13] = 9.075332653468799
$arr[0,14] = 12.54397119453234
$arr[0,15] = -12.79800240606736
$arr[0,16] = -3.829682313083997
$arr[0,17] = -0.1088360584230113
$arr[0,18] = 7.893833719290445
$arr[0,19] = 9.347706739494658
$arr[0,20] = 10.48938412452083
$arr[0,21] = -20.01345706802031
$arr[0,22] = -7.103493917402907
$arr[0,23] = -0.6139735977514333
$arr[0,24] = -18.65772957169785
$arr[0,25] = 10.22202896656103
$arr[0,26] = 6.459959547148602
$arr[0,27] = 3.699109108859769
$arr[0,28] = -0.8523620783265802
$arr[0,29] = -0.8710059723279313
$arr[0,30] = 0.2324312823225227

$arr[1,0] = -2.978358019637952
$arr[1,1] = 1.46110522798605
$arr[1,2] = 1.76488022812836
$arr[1,3] = -0.1819887215705132
$arr[1,4] = -0.06490794078507493
$arr[1,5] = -0.0007307741227983665
$arr[1,6] = -7.557823177607116
$arr[1,7] = 4.290037286496975
$arr[1,8] = 5.37823968090483
$arr[1,9] = -1.613803115652099
$arr[1,10] = -0.4906743893444865
$arr[1,11] = -0.005976284788717225
$arr[1,12] = -4.882783070383476
$arr[1,13] = 9.075332653477679
$arr[1,14] = 12.5439711945501
$arr[1,15] = -12.79800240606514
$arr[1,16] = -3.829682313083441
$arr[1,17] = -0.108836058422994
$arr[1,18] = 7.89383371925492
$arr[1,19] = 9.34770673951242
$arr[1,20] = 10.48938412452971
$arr[1,21] = -20.01345706801143
$arr[1,22] = -7.103493917400686
$arr[1,23] = -0.6139735977512946
$arr[1,24] = -18.65772957173337
$arr[1,25] = 10.22202896656547
$arr[1,26] = 6.459959547157483
$arr[1,27] = 3.699109108877533
$arr[1,28] = -0.8523620783176982
$arr[1,29] = -0.871005972326821
$arr[1,30] = 0.2324312823226364

$arr[2,0] = -3.225733044734957
$arr[2,1] = 1.412074164503649
$arr[2,2] = 1.646627961855945
$arr[2,3] = 0.141646390334282
$arr[2,4] = 0.02510842766352472
$arr[2,5] = 0.0002761003779703362
$arr[2,6] = -8.722718349465147
$arr[2,7] = 3.609298674478238
$arr[2,8] = 4.129622177110615
$arr[2,9] = 0.7966373946401536
$arr[2,10] = 0.1849327619165354
$arr[2,11] = 0.00222734131918904
$arr[2,12] = -16.11529125383388
$arr[2,13] = 4.808850115160541
$arr[2,14] = 4.989208589553469
$arr[2,15] = 4.902080449535722
$arr[2,16] = 1.376955804706298
$arr[2,17] = 0.038196294851406
$arr[2,18] = -9.456536497396705
$arr[2,19] = 0.8474596023921247
$arr[2,20] = 1.300286253585226
$arr[2,21] = 5.250186914725215
$arr[2,22] = 1.879489640070679
$arr[2,23] = 0.1791140865794251
$arr[2,24] = 22.69675951681078
$arr[2,25] = -4.151327021478534
$arr[2,26] = -3.811936822223501
$arr[2,27] = -11.21234432104836
$arr[2,28] = -3.387329281664682
$arr[2,29] = -0.1338220704739709
$arr[2,30] = -0.06083274288090479

$arr[3,0] = 0.2473750310259056
$arr[3,1] = 0.04903106082133424
$arr[3,2] = 0.1182522631433644
$arr[3,3] = -0.3236351120329273
$arr[3,4] = -0.09001636845753576
$arr[3,5] = -0.001006874500860886
$arr[3,6] = 1.164895191223389
$arr[3,7] = 0.6807386034152381
$arr[3,8] = 1.248617493715161
$arr[3,9] = -2.410440510896091
$arr[3,10] = -0.6756071513425111
$arr[3,11] = -0.008203626108851017
$arr[3,12] = 11.2325082463705
$arr[3,13] = 4.266482512845723
$arr[3,14] = 7.554762575856923
$arr[3,15] = -17.70008286225281
$arr[3,16] = -5.206638119442591
$arr[3,17] = -0.1470323533152658
$arr[3,18] = 17.3503703203209
$arr[3,19] = 8.500247105664979
$arr[3,20] = 9.189097835767621
$arr[3,21] = -25.26364401080313
$arr[3,22] = -8.982983565993445
$arr[3,23] = -0.7930876848104748
$arr[3,24] = -41.35448896743232
$arr[3,25] = 14.373355972674
$arr[3,26] = 10.27189634464736
$arr[3,27] = 14.91145337324258
$arr[3,28] = 2.534967181855375
$arr[3,29] = -0.7371839046147971
$arr[3,30] = 0.2932640248045573

$arr[4,0] = 0.2473750310259056
$arr[4,1] = 0.04903106082133452
$arr[4,2] = 0.118252263143365
$arr[4,3] = -0.3236351120331354
$arr[4,4] = -0.0900163684575878
$arr[4,5] = -0.001006874500861462
$arr[4,6] = 1.164895191232268
$arr[4,7] = 0.6807386034152392
$arr[4,8] = 1.248617493710722
$arr[4,9] = -2.410440510897756
$arr[4,10] = -0.6756071513429622
$arr[4,11] = -0.008203626108856167
$arr[4,12] = 11.23250824638827
$arr[4,13] = 4.266482512845723
$arr[4,14] = 7.554762575856925
$arr[4,15] = -17.70008286226392
$arr[4,16] = -5.206638119446477
$arr[4,17] = -0.1470323533153526
$arr[4,18] = 17.35037032039195
$arr[4,19] = 8.500247105664981
$arr[4,20] = 9.189097835749859
$arr[4,21] = -25.26364401083865
$arr[4,22] = -8.982983566006768
$arr[4,23] = -0.7930876848111686
$arr[4,24] = -41.3544889673968
$arr[4,25] = 14.373355972674
$arr[4,26] = 10.27189634463848
$arr[4,27] = 14.91145337320705
$arr[4,28] = 2.534967181837611
$arr[4,29] = -0.7371839046170175
$arr[4,30] = 0.293264024804671

$arr[5,0] = 20936.30056220723
$arr[5,1] = -11003.16588186938
$arr[5,2] = -13544.50428037218
$arr[5,3] = 2735.292079034157
$arr[5,4] = 866.3505052805672
$arr[5,5] = 9.727015739686268
$arr[5,6] = 46342.32248609276
$arr[5,7] = -32639.76689863782
$arr[5,8] = -42784.64650558322
$arr[5,9] = 22447.62345273199
$arr[5,10] = 6554.754566386581
$arr[5,11] = 79.71289891378049
$arr[5,12] = -59565.2963426975
$arr[5,13] = -66872.45674471387
$arr[5,14] = -104177.1889922979
$arr[5,15] = 176861.7041767038
$arr[5,16] = 52277.57973601871
$arr[5,17] = 1475.658166386939
$arr[5,18] = -250591.0742128456
$arr[5,19] = -72456.86879692996
$arr[5,20] = -80563.50791136983
$arr[5,21] = 292878.129953269
$arr[5,22] = 102184.3366030503
$arr[5,23] = 8548.984363075127
$arr[5,24] = 131434.7225201367
$arr[5,25] = -115772.4580324689
$arr[5,26] = -57169.16717272591
$arr[5,27] = -1425.941830326752
$arr[5,28] = 29140.48394512023
$arr[5,29] = 13792.36056686241
$arr[5,30] = -2487.027775987983

$arr[6,0] = 20936.053187176
$arr[6,1] = -11003.21491294887
$arr[6,2] = -13544.62253264513
$arr[6,3] = 2735.615714145775
$arr[6,4] = 866.4405216489994
$arr[6,5] = 9.728022614186997
$arr[6,6] = 46341.15759093119
$arr[6,7] = -32640.44763723585
$arr[6,8] = -42785.8951230825
$arr[6,9] = 22450.03389324309
$arr[6,10] = 6555.430173537885
$arr[6,11] = 79.72110253988677
$arr[6,12] = -59576.52885086241
$arr[6,13] = -66876.72322729877
$arr[6,14] = -104184.7437548046
$arr[6,15] = 176879.4042595413
$arr[6,16] = 52282.78637413495
$arr[6,17] = 1475.805198740255
$arr[6,18] = -250608.424582718
$arr[6,19] = -72465.36904396492
$arr[6,20] = -80572.69700921196
$arr[6,21] = 292903.3935971843
$arr[6,22] = 102193.3195865689
$arr[6,23] = 8549.777450761376
$arr[6,24] = 131476.0770092095
$arr[6,25] = -115786.8313885381
$arr[6,26] = -57179.43906914642
$arr[6,27] = -1440.853283721933
$arr[6,28] = 29137.94897784084
$arr[6,29] = 13793.09775077197
$arr[6,30] = -2487.321040011942

$arr[7,0] = -3.225733046709637
$arr[7,1] = 1.412074165389979
$arr[7,2] = 1.646627962898229
$arr[7,3] = 0.1416463903768797
$arr[7,4] = 0.02510842766647578
$arr[7,5] = 0.0002761003780007607
$arr[7,6] = -8.722718355916651
$arr[7,7] = 3.609298677343222
$arr[7,8] = 4.129622180468449
$arr[7,9] = 0.7966373948405974
$arr[7,10] = 0.1849327619434922
$arr[7,11] = 0.002227341319501418
$arr[7,12] = -16.115291274784
$arr[7,13] = 4.808850123643972
$arr[7,14] = 4.989208599260214
$arr[7,15] = 4.902080451744999
$arr[7,16] = 1.376955805255432
$arr[7,17] = 0.03819629486497573
$arr[7,18] = -9.456536531927652
$arr[7,19] = 0.8474596128631194
$arr[7,20] = 1.300286265304007
$arr[7,21] = 5.250186924062131
$arr[7,22] = 1.879489642907289
$arr[7,23] = 0.1791140867390197
$arr[7,24] = 22.69675947646948
$arr[7,25] = -4.151327016356072
$arr[7,26] = -3.811936813987968
$arr[7,27] = -11.21234430216537
$arr[7,28] = -3.387329274505847
$arr[7,29] = -0.133822069554415
$arr[7,30] = -0.06083274274800488

$ws.Range("B2:AF9").Value = $arr
